$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1328
$ws1.Range("F3").Value = 1218
$ws1.Range("F4").Value = 14606
$ws1.Range("F5").Value = 17681
$ws1.Range("F6").Value = 145
$ws1.Range("F7").Value = 71
$ws1.Range("F9").Value = 216
$ws1.Range("F15").Value = 43
$ws1.Range("F16").Value = 43
$ws1.Range("F17").Value = 147
$ws1.Range("F19").Value = 1333
$ws1.Range("F20").Value = 148
$ws1.Range("F24").Value = 7277
$ws1.Range("F26").Value = 2
$ws1.Range("F27").Value = 38
$ws1.Range("F28").Value = 1174
$ws1.Range("F29").Value = 12
$ws1.Range("F30").Value = 5865
$ws1.Range("F33").Value = 141
$ws1.Range("F35").Value = 226
$ws1.Range("F36").Value = 5093
$ws1.Range("F38").Value = 32

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1328
$ws4.Range("F3").Value = 1218
$ws4.Range("F4").Value = 14606
$ws4.Range("F5").Value = 17681
$ws4.Range("F6").Value = 145
$ws4.Range("F7").Value = 71
$ws4.Range("F9").Value = 216
$ws4.Range("F15").Value = 43
$ws4.Range("F16").Value = 43
$ws4.Range("F17").Value = 147
$ws4.Range("F19").Value = 1333
$ws4.Range("F20").Value = 148
$ws4.Range("F25").Value = 7277
$ws4.Range("F27").Value = 2
$ws4.Range("F28").Value = 38
$ws4.Range("F29").Value = 1174
$ws4.Range("F30").Value = 12
$ws4.Range("F32").Value = 5865
$ws4.Range("F35").Value = 141
$ws4.Range("F37").Value = 226
$ws4.Range("F38").Value = 5093
$ws4.Range("F40").Value = 32

$wb.Save()
